# Actualización automática del tracker
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: write a value into a cell as plain text, without leaving a
# persistent custom number-format style behind (keeps default "Normal" style).
function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# --- Fill in results for already-played matches (rows 31-33) ---

# Row 31: Michael Zheng vs Daniel Masur -> Fallo
Set-TextCell 31 7 "Fallo"
$ws.Cells.Item(31, 8).Value = -1

# Row 32: Kris van Wyk vs Liam Broady -> Fallo
Set-TextCell 32 7 "Fallo"
$ws.Cells.Item(32, 8).Value = -1

# Row 33: Martin Damm Jr vs Samir Banerjee -> Acierto
Set-TextCell 33 7 "Acierto"
$ws.Cells.Item(33, 8).Value = 1

# --- Append new upcoming matches (rows 34-35), results pending ---

# Row 34
$ws.Cells.Item(34, 1).Value = 14728599
Set-TextCell 34 2 "2025-09-21"
$ws.Cells.Item(34, 3).Value = "Taro Daniel"
$ws.Cells.Item(34, 4).Value = "Alexander Shevchenko"
$ws.Cells.Item(34, 5).Value = "Gana Alexander Shevchenko"
$ws.Cells.Item(34, 6).Value = 2
$ws.Cells.Item(34, 7).Font.Bold = $false
$ws.Cells.Item(34, 8).Font.Bold = $false

# Row 35
$ws.Cells.Item(35, 1).Value = 14729132
Set-TextCell 35 2 "2025-09-21"
$ws.Cells.Item(35, 3).Value = "Kaichi Uchida"
$ws.Cells.Item(35, 4).Value = "Hikaru Shiraishi"
$ws.Cells.Item(35, 5).Value = "Gana Hikaru Shiraishi"
$ws.Cells.Item(35, 6).Value = 2.25
$ws.Cells.Item(35, 7).Font.Bold = $false
$ws.Cells.Item(35, 8).Font.Bold = $false
